$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 2785.8667
$ws.Range("I80").Value = 1564.4375
$ws.Range("J80").Value = 4181.7856
$ws.Range("K80").Value = 4693.3125
$ws.Range("L80").Value = 12545.3568
$ws.Range("M80").Value = -3695.3125
$ws.Range("N80").Value = -14541.3568
$ws.Range("H83").Value = 2785.8667
$ws.Range("I83").Value = 1564.4375
$ws.Range("J83").Value = 4181.7856
$ws.Range("K83").Value = 14079.9375
$ws.Range("L83").Value = 37636.0704
$ws.Range("M83").Value = -9087.9375
$ws.Range("N83").Value = -47620.0704
$ws.Range("H98").Value = 188802.81
$ws.Range("I98").Value = 1023.7241
$ws.Range("K98").Value = 1023.7241
$ws.Range("M98").Value = 474.2759
$ws.Range("H122").Value = 188802.81
$ws.Range("I122").Value = 1023.7241
$ws.Range("K122").Value = 3071.1723
$ws.Range("M122").Value = -621.1723000000002
$ws.Range("H137").Value = 12989701
$ws.Range("I137").Value = 47620764
$ws.Range("J137").Value = 3052.375
$ws.Range("K137").Value = 142862292
$ws.Range("L137").Value = 9157.125
$ws.Range("M137").Value = -142859742
$ws.Range("N137").Value = -14257.125
$ws.Range("H138").Value = 2830.8064
$ws.Range("J138").Value = 3106.589
$ws.Range("L138").Value = 9319.767
$ws.Range("N138").Value = -19599.767

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4565.533
$ws.Range("I2").Value = 795.03705
$ws.Range("K2").Value = 795.03705
$ws.Range("M2").Value = -682.03705
$ws.Range("H32").Value = 4048.027
$ws.Range("I32").Value = 3066.232
$ws.Range("K32").Value = 3066.232
$ws.Range("M32").Value = -2779.232
$ws.Range("H45").Value = 111112860
$ws.Range("I45").Value = 111112860
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 111112860
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -111112483
$ws.Range("N45").ClearContents()
$ws.Range("H61").Value = 7549.5557
$ws.Range("I61").Value = 4954.5713
$ws.Range("K61").Value = 4954.5713
$ws.Range("M61").Value = -4742.5713
$ws.Range("H74").Value = 16669956
$ws.Range("J74").Value = 4736
$ws.Range("L74").Value = 4736
$ws.Range("N74").Value = -6484
$ws.Range("H77").Value = 16669956
$ws.Range("J77").Value = 4736
$ws.Range("L77").Value = 23680
$ws.Range("N77").Value = -32416
$ws.Range("H110").Value = 6723.25
$ws.Range("I110").Value = 6036.8
$ws.Range("K110").Value = 6036.8
$ws.Range("M110").Value = -3991.8
$ws.Range("H116").Value = 4565.533
$ws.Range("I116").Value = 795.03705
$ws.Range("K116").Value = 795.03705
$ws.Range("M116").Value = 1498.96295
$ws.Range("H136").Value = 7549.5557
$ws.Range("I136").Value = 4954.5713
$ws.Range("K136").Value = 14863.7139
$ws.Range("M136").Value = -12313.7139
$ws.Range("H139").Value = 68887.78
$ws.Range("I139").Value = 68589
$ws.Range("J139").Value = 68925.125
$ws.Range("K139").Value = 68589
$ws.Range("L139").Value = 68925.125
$ws.Range("M139").Value = -63449
$ws.Range("N139").Value = -79205.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4565.533
$ws.Range("I3").Value = 795.03705
$ws.Range("K3").Value = 795.03705
$ws.Range("M3").Value = -681.03705
$ws.Range("H7").Value = 4267.75
$ws.Range("I7").Value = 5023.6665
$ws.Range("J7").Value = 2000
$ws.Range("K7").Value = 5023.6665
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = -4910.6665
$ws.Range("N7").Value = -2226
$ws.Range("H20").Value = 4601.44
$ws.Range("I20").Value = 4303.4707
$ws.Range("J20").Value = 5234.625
$ws.Range("K20").Value = 4303.4707
$ws.Range("L20").Value = 5234.625
$ws.Range("M20").Value = -4056.4707
$ws.Range("N20").Value = -5728.625
$ws.Range("H86").Value = 3451.32
$ws.Range("I86").Value = 3514.4614
$ws.Range("J86").Value = 3382.9167
$ws.Range("K86").Value = 3514.4614
$ws.Range("L86").Value = 3382.9167
$ws.Range("M86").Value = -2391.4614
$ws.Range("N86").Value = -5628.9167
$ws.Range("H89").Value = 3451.32
$ws.Range("I89").Value = 3514.4614
$ws.Range("J89").Value = 3382.9167
$ws.Range("K89").Value = 17572.307
$ws.Range("L89").Value = 16914.5835
$ws.Range("M89").Value = -11956.307
$ws.Range("N89").Value = -28146.5835
$ws.Range("H138").Value = 65116.5
$ws.Range("J138").Value = 65116.5
$ws.Range("L138").Value = 65116.5
$ws.Range("N138").Value = -75396.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1252.375
$ws.Range("I16").Value = 632
$ws.Range("K16").Value = 632
$ws.Range("M16").Value = -345
$ws.Range("H31").Value = 27692.691
$ws.Range("I31").Value = 1926.2693
$ws.Range("J31").Value = 69563.125
$ws.Range("K31").Value = 1926.2693
$ws.Range("L31").Value = 69563.125
$ws.Range("M31").Value = -1631.2693
$ws.Range("N31").Value = -70153.125
$ws.Range("H34").Value = 27692.691
$ws.Range("I34").Value = 1926.2693
$ws.Range("J34").Value = 69563.125
$ws.Range("K34").Value = 1926.2693
$ws.Range("L34").Value = 69563.125
$ws.Range("M34").Value = -1724.2693
$ws.Range("N34").Value = -69967.125
$ws.Range("H99").Value = 2613.8572
$ws.Range("I99").Value = 2150
$ws.Range("K99").Value = 2150
$ws.Range("M99").Value = -652
$ws.Range("H113").Value = 1252.375
$ws.Range("I113").Value = 632
$ws.Range("K113").Value = 632
$ws.Range("M113").Value = 1538
$ws.Range("H126").Value = 2613.8572
$ws.Range("I126").Value = 2150
$ws.Range("K126").Value = 6450
$ws.Range("M126").Value = -3980

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 62898.625
$ws.Range("I7").Value = 498
$ws.Range("J7").Value = 166899.67
$ws.Range("K7").Value = 1494
$ws.Range("L7").Value = 500699.01
$ws.Range("M7").Value = -1382
$ws.Range("N7").Value = -500923.01

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12636.728
$ws.Range("J70").Value = 15201.4
$ws.Range("L70").Value = 15201.4
$ws.Range("N70").Value = -15741.4
$ws.Range("H73").Value = 12636.728
$ws.Range("J73").Value = 15201.4
$ws.Range("L73").Value = 15201.4
$ws.Range("N73").Value = -17073.4
$ws.Range("H80").Value = 459163
$ws.Range("I80").Value = 558420.75
$ws.Range("K80").Value = 558420.75
$ws.Range("M80").Value = -557422.75
$ws.Range("H83").Value = 459163
$ws.Range("I83").Value = 558420.75
$ws.Range("K83").Value = 2792103.75
$ws.Range("M83").Value = -2787111.75
$ws.Range("H107").Value = 824.9091
$ws.Range("I107").Value = 915.625
$ws.Range("J107").Value = 583
$ws.Range("K107").Value = 915.625
$ws.Range("L107").Value = 583
$ws.Range("M107").Value = 1004.375
$ws.Range("N107").Value = -4423

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2670.4263
$ws.Range("I136").Value = 1706.1915
$ws.Range("K136").Value = 5118.5745
$ws.Range("M136").Value = -2568.5745

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 1000500
$ws.Range("J8").Value = 1000500
$ws.Range("L8").Value = 1000500
$ws.Range("N8").Value = -1000780
$ws.Range("H81").Value = 4251.4546
$ws.Range("I81").Value = 2775
$ws.Range("J81").Value = 8188.6665
$ws.Range("K81").Value = 5550
$ws.Range("L81").Value = 16377.333
$ws.Range("M81").Value = -4489
$ws.Range("N81").Value = -18499.333
$ws.Range("H84").Value = 4251.4546
$ws.Range("I84").Value = 2775
$ws.Range("J84").Value = 8188.6665
$ws.Range("K84").Value = 27750
$ws.Range("L84").Value = 81886.66500000001
$ws.Range("M84").Value = -22446
$ws.Range("N84").Value = -92494.66500000001
$ws.Range("H107").Value = 1019.6923
$ws.Range("I107").Value = 854.6667
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 2564.0001
$ws.Range("L107").Value = 9000
$ws.Range("M107").Value = -644.0001000000002
$ws.Range("N107").Value = -12840
$ws.Range("H136").Value = 3284.383
$ws.Range("I136").Value = 2126.0571
$ws.Range("J136").Value = 6662.8335
$ws.Range("K136").Value = 6378.1713
$ws.Range("L136").Value = 19988.5005
$ws.Range("M136").Value = -3828.1713
$ws.Range("N136").Value = -25088.5005
